$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.162809
$ws.Range("H2").Value = 6.488427000000001
$ws.Range("I2").Value = 0.06755089002018773
$ws.Range("J2").Value = 0.06755089002018773
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.103724333333333
$ws.Range("N2").Value = 3.311173
$ws.Range("O2").Value = 0.01870879385910814
$ws.Range("P2").Value = 0.01870879385910814
$ws.Range("Q2").Value = 2.387144921652334
$ws.Range("R2").Value = 21.484304294871
$ws.Range("S2").Value = 0.001263795676386978
$ws.Range("T2").Value = 0.001263795676386978
$ws.Range("G3").Value = 2.162809
$ws.Range("H3").Value = 6.488427000000001
$ws.Range("I3").Value = 0.06755089002018773
$ws.Range("J3").Value = 0.06755089002018773
$ws.Range("O3").Value = 0.1603368629650925
$ws.Range("P3").Value = 0.1603368629650925
$ws.Range("Q3").Value = 20.458150913585
$ws.Range("R3").Value = 184.123358222265
$ws.Range("S3").Value = 0.01083089779633688
$ws.Range("T3").Value = 0.01083089779633688
$ws.Range("G4").Value = 2.162809
$ws.Range("H4").Value = 6.488427000000001
$ws.Range("I4").Value = 0.06755089002018773
$ws.Range("J4").Value = 0.06755089002018773
$ws.Range("M4").Value = 47.61312599999999
$ws.Range("N4").Value = 142.839378
$ws.Range("O4").Value = 0.80707123365805
$ws.Range("P4").Value = 0.80707123365805
$ws.Range("Q4").Value = 102.978097430934
$ws.Range("R4").Value = 926.802876878406
$ws.Range("S4").Value = 0.05451838014329217
$ws.Range("T4").Value = 0.05451838014329217
$ws.Range("G5").Value = 2.162809
$ws.Range("H5").Value = 6.488427000000001
$ws.Range("I5").Value = 0.06755089002018773
$ws.Range("J5").Value = 0.06755089002018773
$ws.Range("M5").Value = 0.8190333333333334
$ws.Range("N5").Value = 2.4571
$ws.Range("O5").Value = 0.01388310951774934
$ws.Range("P5").Value = 0.01388310951774934
$ws.Range("Q5").Value = 1.771412664633334
$ws.Range("R5").Value = 15.9427139817
$ws.Range("S5").Value = 0.0009378164041717069
$ws.Range("T5").Value = 0.0009378164041717069
$ws.Range("I6").Value = 0.5628021396814664
$ws.Range("J6").Value = 0.5628021396814664
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.103724333333333
$ws.Range("N6").Value = 3.311173
$ws.Range("O6").Value = 0.01870879385910814
$ws.Range("P6").Value = 0.01870879385910814
$ws.Range("Q6").Value = 19.88856503939733
$ws.Range("R6").Value = 178.997085354576
$ws.Range("S6").Value = 0.01052934921476554
$ws.Range("T6").Value = 0.01052934921476554
$ws.Range("I7").Value = 0.5628021396814664
$ws.Range("J7").Value = 0.5628021396814664
$ws.Range("O7").Value = 0.1603368629650925
$ws.Range("P7").Value = 0.1603368629650925
$ws.Range("S7").Value = 0.09023792954656815
$ws.Range("T7").Value = 0.09023792954656815
$ws.Range("I8").Value = 0.5628021396814664
$ws.Range("J8").Value = 0.5628021396814664
$ws.Range("M8").Value = 47.61312599999999
$ws.Range("N8").Value = 142.839378
$ws.Range("O8").Value = 0.80707123365805
$ws.Range("P8").Value = 0.80707123365805
$ws.Range("Q8").Value = 857.964914409504
$ws.Range("R8").Value = 7721.684229685535
$ws.Range("S8").Value = 0.4542214171781113
$ws.Range("T8").Value = 0.4542214171781113
$ws.Range("I9").Value = 0.5628021396814664
$ws.Range("J9").Value = 0.5628021396814664
$ws.Range("M9").Value = 0.8190333333333334
$ws.Range("N9").Value = 2.4571
$ws.Range("O9").Value = 0.01388310951774934
$ws.Range("P9").Value = 0.01388310951774934
$ws.Range("Q9").Value = 14.75857442613334
$ws.Range("R9").Value = 132.8271698352
$ws.Range("S9").Value = 0.007813443742021458
$ws.Range("T9").Value = 0.007813443742021458
$ws.Range("G10").Value = 4.650307000000001
$ws.Range("H10").Value = 13.950921
$ws.Range("I10").Value = 0.1452427730405732
$ws.Range("J10").Value = 0.1452427730405732
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.103724333333333
$ws.Range("N10").Value = 3.311173
$ws.Range("O10").Value = 0.01870879385910814
$ws.Range("P10").Value = 0.01870879385910814
$ws.Range("Q10").Value = 5.132656993370334
$ws.Range("R10").Value = 46.19391294033301
$ws.Range("S10").Value = 0.002717317100341314
$ws.Range("T10").Value = 0.002717317100341314
$ws.Range("G11").Value = 4.650307000000001
$ws.Range("H11").Value = 13.950921
$ws.Range("I11").Value = 0.1452427730405732
$ws.Range("J11").Value = 0.1452427730405732
$ws.Range("O11").Value = 0.1603368629650925
$ws.Range("P11").Value = 0.1603368629650925
$ws.Range("Q11").Value = 43.98755618295501
$ws.Range("R11").Value = 395.888005646595
$ws.Range("S11").Value = 0.02328777059767643
$ws.Range("T11").Value = 0.02328777059767643
$ws.Range("G12").Value = 4.650307000000001
$ws.Range("H12").Value = 13.950921
$ws.Range("I12").Value = 0.1452427730405732
$ws.Range("J12").Value = 0.1452427730405732
$ws.Range("M12").Value = 47.61312599999999
$ws.Range("N12").Value = 142.839378
$ws.Range("O12").Value = 0.80707123365805
$ws.Range("P12").Value = 0.80707123365805
$ws.Range("Q12").Value = 221.415653129682
$ws.Range("R12").Value = 1992.740878167138
$ws.Range("S12").Value = 0.1172212640177716
$ws.Range("T12").Value = 0.1172212640177716
$ws.Range("G13").Value = 4.650307000000001
$ws.Range("H13").Value = 13.950921
$ws.Range("I13").Value = 0.1452427730405732
$ws.Range("J13").Value = 0.1452427730405732
$ws.Range("M13").Value = 0.8190333333333334
$ws.Range("N13").Value = 2.4571
$ws.Range("O13").Value = 0.01388310951774934
$ws.Range("P13").Value = 0.01388310951774934
$ws.Range("Q13").Value = 3.808756443233334
$ws.Range("R13").Value = 34.2788079891
$ws.Range("S13").Value = 0.002016421324783889
$ws.Range("T13").Value = 0.002016421324783889
$ws.Range("G14").Value = 7.184856000000001
$ws.Range("H14").Value = 21.554568
$ws.Range("I14").Value = 0.2244041972577726
$ws.Range("J14").Value = 0.2244041972577726
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 1.103724333333333
$ws.Range("N14").Value = 3.311173
$ws.Range("O14").Value = 0.01870879385910814
$ws.Range("P14").Value = 0.01870879385910814
$ws.Range("Q14").Value = 7.930100398696001
$ws.Range("R14").Value = 71.37090358826401
$ws.Range("S14").Value = 0.004198331867614309
$ws.Range("T14").Value = 0.004198331867614309
$ws.Range("G15").Value = 7.184856000000001
$ws.Range("H15").Value = 21.554568
$ws.Range("I15").Value = 0.2244041972577726
$ws.Range("J15").Value = 0.2244041972577726
$ws.Range("O15").Value = 0.1603368629650925
$ws.Range("P15").Value = 0.1603368629650925
$ws.Range("Q15").Value = 67.96201991964001
$ws.Range("R15").Value = 611.6581792767601
$ws.Range("S15").Value = 0.03598026502451108
$ws.Range("T15").Value = 0.03598026502451108
$ws.Range("G16").Value = 7.184856000000001
$ws.Range("H16").Value = 21.554568
$ws.Range("I16").Value = 0.2244041972577726
$ws.Range("J16").Value = 0.2244041972577726
$ws.Range("M16").Value = 47.61312599999999
$ws.Range("N16").Value = 142.839378
$ws.Range("O16").Value = 0.80707123365805
$ws.Range("P16").Value = 0.80707123365805
$ws.Range("Q16").Value = 342.093454019856
$ws.Range("R16").Value = 3078.841086178704
$ws.Range("S16").Value = 0.1811101723188749
$ws.Range("T16").Value = 0.1811101723188749
$ws.Range("G17").Value = 7.184856000000001
$ws.Range("H17").Value = 21.554568
$ws.Range("I17").Value = 0.2244041972577726
$ws.Range("J17").Value = 0.2244041972577726
$ws.Range("M17").Value = 0.8190333333333334
$ws.Range("N17").Value = 2.4571
$ws.Range("O17").Value = 0.01388310951774934
$ws.Range("P17").Value = 0.01388310951774934
$ws.Range("Q17").Value = 5.884636559200001
$ws.Range("R17").Value = 52.96172903280001
$ws.Range("S17").Value = 0.003115428046772282
$ws.Range("T17").Value = 0.003115428046772282
